# Fill in the previously-empty "User Interface" (28.3), "Velocity" (28.4),
# and "Orbital Element" (28.5) class-diagram tables with the class name,
# fields, and operations text, matching the Consolas-font "Answer" style
# already established by the other class tables in the document.

$d = $word.ActiveDocument

function Fill-Cell($table, $row, $col, $text) {
    $cell = $table.Cell($row, $col)
    $range = $cell.Range
    $range.Text = $text
    $cell.Range.Font.Name = "Consolas"
}

# --- Problem 28.3: User Interface ---
$t3 = $d.Tables.Item(3)
Fill-Cell $t3 1 1 "Orbital Simulator"
Fill-Cell $t3 2 1 "objectList"
Fill-Cell $t3 3 1 "addObject`rsetObjectPosition`rsetObjectState`rdisplayStatus"

# --- Problem 28.4: Velocity ---
$t4 = $d.Tables.Item(4)
Fill-Cell $t4 1 1 "Velocity"
Fill-Cell $t4 2 1 "Speed`rDirection"
Fill-Cell $t4 3 1 "increaseVelocity`rmultiplyVelocity`rreverseVelocity`rdivideVelocity`roperator +`roperator –`roperator *`roperator /`roperator = "

# --- Problem 28.5: Orbital Element ---
$t5 = $d.Tables.Item(5)
Fill-Cell $t5 1 1 "Element"
Fill-Cell $t5 2 1 "Position position`rVelocity velocity`rstatus"
Fill-Cell $t5 3 1 "hasCollided`rdrawElement`rmoveElement`risAlive"

Write-Host "Done filling tables"
